# LAB4: add test for method ureadCar(Integer id)
# Update the stale "car{id}" URL in C5 to the correct "car/{id}" URL
# (matches the URL already used in C6/C7), and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "http://localhost:8080/api/car/{id}"

$ws.Range("C17").Select()
